$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("EU28+CH") was removed entirely; columns F:L shift left to become E:K.
$ws.Columns("E:E").Delete()

# Corrected counts ("hatching" error fix) in column B plus one corrected value in J10.
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("B10").ClearContents()
$ws.Range("J10").Value = 2
$ws.Range("B12").ClearContents()
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 1
